$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.060.11'
$ws.Range('E2').Value = '  +5.67%  '
$ws.Range('D3').Value = '3.603.22'
$ws.Range('E3').Value = '  +5.46%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '190.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.646'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.20%  '
$ws.Range('D8').Value = '3.592.63'
$ws.Range('E8').Value = '  +5.41%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.661'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000291'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.76%  '
$ws.Range('D15').Value = '4.171.33'
$ws.Range('E15').Value = '  +5.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.56%  '
$ws.Range('D17').Value = '3.596.49'
$ws.Range('E17').Value = '  +5.38%  '
$ws.Range('D18').Value = '69.993.28'
$ws.Range('E18').Value = '  +5.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.11%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('E21').Value = '  +4.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '494.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +21.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.55%  '
$ws.Range('E25').Value = '  +8.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '90.69'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('E27').Value = '  +5.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.69'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '618.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.33%  '
$ws.Range('E34').Value = '  +7.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '65.06'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.03%  '
$ws.Range('D36').Value = '0.0₃0823'
$ws.Range('E36').Value = '  +8.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.19'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.15%  '
$ws.Range('E38').Value = '  +5.51%  '
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.64'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('D42').Value = '3.315.97'
$ws.Range('E42').Value = '  +5.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.08'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.25%  '
$ws.Range('E44').Value = '  +5.20%  '
$ws.Range('E45').Value = '  +6.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.14'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.33%  '
$ws.Range('E49').Value = '  -2.45%  '
$ws.Range('E50').Value = '  +5.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.10%  '
